$wb = $excel.ActiveWorkbook

# --- Sheet "Top 10 empresas del momento" ---
$ws3 = $wb.Worksheets.Item("Top 10 empresas del momento")

$ws3.Range("A4").Value  = "Microsoft"
$ws3.Range("A5").Value  = "Amazon"
$ws3.Range("A6").Value  = "Nvidia"
$ws3.Range("A7").Value  = "Zoom"
$ws3.Range("A8").Value  = "Paypal"
$ws3.Range("A9").Value  = "Netflix"
$ws3.Range("A10").Value = "Facebook"
$ws3.Range("A11").Value = "Use Stop Loss and Take Profit to manage your risks"

# --- Sheet "Top 10 acciones bajas" ---
$ws4 = $wb.Worksheets.Item("Top 10 acciones bajas")

$ws4.Range("A2").Value  = "Sundial Growers"
$ws4.Range("A3").Value  = "Zomedica"
$ws4.Range("A4").Value  = "Check-Cap"
$ws4.Range("A5").Value  = "Castor Maritime"
$ws4.Range("A6").Value  = "Cinedigm"
$ws4.Range("A7").Value  = "Atossa Therapeutics"
$ws4.Range("A8").Value  = "Ocugen"
$ws4.Range("A9").Value  = "Fourth Wave Energy"
$ws4.Range("A10").Value = "AbCellera Biologics"
$ws4.Range("A11").Value = "Clean Energy Fuels"
$ws4.Range("A12").Value = "Invest in a diverse set of assets to spread risk"
$ws4.Range("A13").Value = "Top 10 acciones a bajo costo con potencial "
$ws4.Range("A14").Value = "Sundial Growers"
$ws4.Range("A15").Value = "Zomedica"
$ws4.Range("A16").Value = "Check-Cap"
$ws4.Range("A17").Value = "Castor Maritime"
$ws4.Range("A18").Value = "Cinedigm"
$ws4.Range("A19").Value = "Atossa Therapeutics"
$ws4.Range("A20").Value = "Ocugen"
$ws4.Range("A21").Value = "Fourth Wave Energy"
$ws4.Range("A22").Value = "AbCellera Biologics"
$ws4.Range("A23").Value = "Clean Energy Fuels"
$ws4.Range("A24").Value = "Invest in a diverse set of assets to spread risk"
